# Update the Cthrc1-Fzd3 LR-pairs worksheet with newly computed TPM-based values.
# This includes:
#   1. Renaming the "Resolving-Mac" cluster label to "Inflammatory-Mac".
#   2. Updating numeric statistics across rows 2-13 to reflect the new TPM data
#      (ligand/receptor expression values, specificities, edge weights, etc.).
#   3. Rows that referenced the renamed cluster ("MuSCs" vs the renamed cluster)
#      end up with swapped Target-cluster labels/values between the two rows
#      that had K/L/M/N/O/P/Q/R/S/T figures for MuSCs and the renamed cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated data values (recomputed with new TPM input) to the table body.
# This also renames every occurrence of the "Resolving-Mac" cluster label to
# "Inflammatory-Mac" (cells D4, D8, D12 below), while the cells that used to read
# "Resolving-Mac" because of a reordering (D5, D9, D13) are now set to "MuSCs".
$ws.Range("G2").Value = 1.324075666666667
$ws.Range("H2").Value = 3.972227
$ws.Range("I2").Value = 0.01675578032580584
$ws.Range("J2").Value = 0.01684165790066494
$ws.Range("M2").Value = 0.1579376666666667
$ws.Range("N2").Value = 0.473813
$ws.Range("O2").Value = 0.05467876644486869
$ws.Range("P2").Value = 0.07340983674118848
$ws.Range("Q2").Value = 0.2091214212834444
$ws.Range("R2").Value = 1.882092791551
$ws.Range("S2").Value = 0.0009161853990362631
$ws.Range("T2").Value = 0.00123634335693876
$ws.Range("G3").Value = 1.324075666666667
$ws.Range("H3").Value = 3.972227
$ws.Range("I3").Value = 0.01675578032580584
$ws.Range("J3").Value = 0.01684165790066494
$ws.Range("O3").Value = 0.1724539210166233
$ws.Range("P3").Value = 0.2315307204300726
$ws.Range("Q3").Value = 0.6595578396096667
$ws.Range("R3").Value = 5.936020556487001
$ws.Range("S3").Value = 0.00288960001687841
$ws.Range("T3").Value = 0.003899361186977777
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 1.324075666666667
$ws.Range("H4").Value = 3.972227
$ws.Range("I4").Value = 0.01675578032580584
$ws.Range("J4").Value = 0.01684165790066494
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02135966666666667
$ws.Range("N4").Value = 0.064079
$ws.Range("O4").Value = 0.00739481752299059
$ws.Range("P4").Value = 0.0099280284174107
$ws.Range("Q4").Value = 0.02828181488144444
$ws.Range("R4").Value = 0.254536333933
$ws.Range("S4").Value = 0.00012390593796465
$ws.Range("T4").Value = 0.0001672044582341109
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 1.324075666666667
$ws.Range("H5").Value = 3.972227
$ws.Range("I5").Value = 0.01675578032580584
$ws.Range("J5").Value = 0.01684165790066494
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.21104
$ws.Range("N5").Value = 4.422079999999999
$ws.Range("O5").Value = 0.7654724950155174
$ws.Range("P5").Value = 0.6851314144113283
$ws.Range("Q5").Value = 2.927584262026666
$ws.Range("R5").Value = 17.56550557216
$ws.Range("S5").Value = 0.01282608897192651
$ws.Range("T5").Value = 0.01153874889851429
$ws.Range("I6").Value = 0.9679468703219594
$ws.Range("J6").Value = 0.9729078406975189
$ws.Range("M6").Value = 0.1579376666666667
$ws.Range("N6").Value = 0.473813
$ws.Range("O6").Value = 0.05467876644486869
$ws.Range("P6").Value = 0.07340983674118848
$ws.Range("Q6").Value = 12.08051319083255
$ws.Range("R6").Value = 108.724618717493
$ws.Range("S6").Value = 0.05292614085337601
$ws.Range("T6").Value = 0.07142100574982707
$ws.Range("I7").Value = 0.9679468703219594
$ws.Range("J7").Value = 0.9729078406975189
$ws.Range("O7").Value = 0.1724539210166233
$ws.Range("P7").Value = 0.2315307204300726
$ws.Range("R7").Value = 342.911664398541
$ws.Range("S7").Value = 0.1669262331227909
$ws.Range("T7").Value = 0.2252580532687628
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.9679468703219594
$ws.Range("J8").Value = 0.9729078406975189
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02135966666666667
$ws.Range("N8").Value = 0.064079
$ws.Range("O8").Value = 0.00739481752299059
$ws.Range("P8").Value = 0.0099280284174107
$ws.Range("Q8").Value = 1.633782113946555
$ws.Range("R8").Value = 14.704039025519
$ws.Range("S8").Value = 0.007157790477980725
$ws.Range("T8").Value = 0.00965905668996665
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.9679468703219594
$ws.Range("J9").Value = 0.9729078406975189
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.21104
$ws.Range("N9").Value = 4.422079999999999
$ws.Range("O9").Value = 0.7654724950155174
$ws.Range("P9").Value = 0.6851314144113283
$ws.Range("Q9").Value = 169.1205046218133
$ws.Range("R9").Value = 1014.72302773088
$ws.Range("S9").Value = 0.7409367058678117
$ws.Range("T9").Value = 0.6665697249889624
$ws.Range("G10").Value = 1.2088275
$ws.Range("H10").Value = 2.417655
$ws.Range("I10").Value = 0.0152973493522347
$ws.Range("J10").Value = 0.01025050140181618
$ws.Range("M10").Value = 0.1579376666666667
$ws.Range("N10").Value = 0.473813
$ws.Range("O10").Value = 0.05467876644486869
$ws.Range("P10").Value = 0.07340983674118848
$ws.Range("Q10").Value = 0.1909193947525
$ws.Range("R10").Value = 1.145516368515
$ws.Range("S10").Value = 0.0008364401924564043
$ws.Range("T10").Value = 0.0007524876344226496
$ws.Range("G11").Value = 1.2088275
$ws.Range("H11").Value = 2.417655
$ws.Range("I11").Value = 0.0152973493522347
$ws.Range("J11").Value = 0.01025050140181618
$ws.Range("O11").Value = 0.1724539210166233
$ws.Range("P11").Value = 0.2315307204300726
$ws.Range("Q11").Value = 0.6021496160925001
$ws.Range("R11").Value = 3.612897696555
$ws.Range("S11").Value = 0.002638087876953975
$ws.Range("T11").Value = 0.00237330597433197
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 1.2088275
$ws.Range("H12").Value = 2.417655
$ws.Range("I12").Value = 0.0152973493522347
$ws.Range("J12").Value = 0.01025050140181618
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02135966666666667
$ws.Range("N12").Value = 0.064079
$ws.Range("O12").Value = 0.00739481752299059
$ws.Range("P12").Value = 0.0099280284174107
$ws.Range("Q12").Value = 0.0258201524575
$ws.Range("R12").Value = 0.154920914745
$ws.Range("S12").Value = 0.0001131211070452139
$ws.Range("T12").Value = 0.0001017672692099393
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 1.2088275
$ws.Range("H13").Value = 2.417655
$ws.Range("I13").Value = 0.0152973493522347
$ws.Range("J13").Value = 0.01025050140181618
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.21104
$ws.Range("N13").Value = 4.422079999999999
$ws.Range("O13").Value = 0.7654724950155174
$ws.Range("P13").Value = 0.6851314144113283
$ws.Range("Q13").Value = 2.6727659556
$ws.Range("R13").Value = 10.6910638224
$ws.Range("S13").Value = 0.0117097001757791
$ws.Range("T13").Value = 0.007022940523851625
